# Automatic update of files.
# Applies the diff: updates Taxonsorteringsordning (column B) values,
# and swaps the content of rows 5 and 6 (species got re-matched to the
# correct observation row), moving the "Publik kommentar" note along
# with its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Taxonsorteringsordning (column B) bumps ---------------------
$simpleRows = @(2, 3, 7, 9, 10, 11, 12, 13, 14, 15, 17)
foreach ($r in $simpleRows) {
    $ws.Cells.Item($r, 2).Value = 79244
}
$ws.Cells.Item(16, 2).Value = 91830

# --- Swap the data of row 5 and row 6 (all columns except B) ------------
# NOTE: use .Value2 for reading - .Value on this engine can return a
# wrapper object rather than the underlying scalar.
# Capture row 5 values before overwriting anything.
$row5_A  = $ws.Cells.Item(5, 1).Value2     # A  Id
$row5_E  = $ws.Cells.Item(5, 5).Value2     # E  TaxonId
$row5_F  = $ws.Cells.Item(5, 6).Value2     # F  Artnamn
$row5_G  = $ws.Cells.Item(5, 7).Value2     # G  Vetenskapligt namn
$row5_H  = $ws.Cells.Item(5, 8).Value2     # H  Auktor
$row5_J  = $ws.Cells.Item(5, 10).Value2    # J  Enhet
$row5_Q  = $ws.Cells.Item(5, 17).Value2    # Q  Ost
$row5_R  = $ws.Cells.Item(5, 18).Value2    # R  Nord
$row5_AC = $ws.Cells.Item(5, 29).Value2    # AC Publik kommentar

$row6_A  = $ws.Cells.Item(6, 1).Value2
$row6_E  = $ws.Cells.Item(6, 5).Value2
$row6_F  = $ws.Cells.Item(6, 6).Value2
$row6_G  = $ws.Cells.Item(6, 7).Value2
$row6_H  = $ws.Cells.Item(6, 8).Value2
$row6_J  = $ws.Cells.Item(6, 10).Value2
$row6_Q  = $ws.Cells.Item(6, 17).Value2
$row6_R  = $ws.Cells.Item(6, 18).Value2
$row6_AC = $ws.Cells.Item(6, 29).Value2

# Write row 5 <- old row 6 values
$ws.Cells.Item(5, 1).Value  = $row6_A
$ws.Cells.Item(5, 5).Value  = $row6_E
$ws.Cells.Item(5, 6).Value  = $row6_F
$ws.Cells.Item(5, 7).Value  = $row6_G
$ws.Cells.Item(5, 8).Value  = $row6_H
$ws.Cells.Item(5, 10).Value = $row6_J
$ws.Cells.Item(5, 17).Value = $row6_Q
$ws.Cells.Item(5, 18).Value = $row6_R
$ws.Cells.Item(5, 29).Value = ""

# Write row 6 <- old row 5 values
$ws.Cells.Item(6, 1).Value  = $row5_A
$ws.Cells.Item(6, 5).Value  = $row5_E
$ws.Cells.Item(6, 6).Value  = $row5_F
$ws.Cells.Item(6, 7).Value  = $row5_G
$ws.Cells.Item(6, 8).Value  = $row5_H
$ws.Cells.Item(6, 10).Value = $row5_J
$ws.Cells.Item(6, 17).Value = $row5_Q
$ws.Cells.Item(6, 18).Value = $row5_R
$ws.Cells.Item(6, 29).Value = $row5_AC

# Row 5's own Taxonsorteringsordning becomes 79244, row 6's becomes 83090
$ws.Cells.Item(5, 2).Value = 79244
$ws.Cells.Item(6, 2).Value = 83090
